$p = $ppt.ActivePresentation

# --- Slide 2: "Acknowledgment " -> "Acknowledgments " ---
$s2 = $p.Slides.Item(2)
$titleShape = $s2.Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Acknowledgments "

# --- Slide 7: Title shape resize/reposition + text update ---
$s7 = $p.Slides.Item(7)
$shmemTitle = $s7.Shapes.Item(2)

# Reposition / resize the title placeholder
$shmemTitle.Left = 27.5
$shmemTitle.Top = 7.5
$shmemTitle.Width = 675.8443
$shmemTitle.Height = 45

$tr = $shmemTitle.TextFrame.TextRange

# "Shmem" (chars 1-5) -> "OpenShmem"
$tr.Characters(1, 5).Text = "OpenShmem"

# After the above change the leading space of the former " Analyzer in " run
# sits at position 10; splitting it into its own run (" ") and shrinking the
# remainder to "Analyzer in " matches the target run layout.
$tr.Characters(11, 12).Text = "Analyzer in "
